# Revert CDRDfRCP parameters to avoid artificially high CES credit prices
#
# - CDRDfRCP!B1 (the "a" parameter) goes from 75 back to 15
# - CDRDfRCP!B2 (the "b" parameter) goes from 0.98 back to 0.9
# - About sheet gains two explanatory note rows (10 & 11) under the
#   existing "Notes:" section
# - Formula-driven values on About!B15:B63 recompute from the new
#   parameters (handled automatically by the workbook's auto-recalc)
# - Selections on both sheets move to reflect where the editor last
#   clicked

$wb = $excel.ActiveWorkbook

$wsCDR = $wb.Worksheets.Item("CDRDfRCP")
$wsAbout = $wb.Worksheets.Item("About")

# --- Core parameter edits (drives the dependent formulas + chart) ---
$wsCDR.Range("B1").Value = 15
$wsCDR.Range("B2").Value = 0.9

# --- New explanatory notes under the "Notes:" section on About ---
$wsAbout.Range("B10").Value = "We avoid having this increase too sharply in the last few percent of requirement to avoid"
$wsAbout.Range("B11").Value = "artificially high CES credit prices in the model."

# --- Restore the on-screen selections to match the saved state ---
# (CDRDfRCP is selected first so that "About" ends up as the active tab,
# matching the workbook's tabSelected flag.)
$wsCDR.Activate()
$wsCDR.Range("B3").Select()

$wsAbout.Activate()
$wsAbout.Range("B12").Select()
